$wb = $excel.ActiveWorkbook

# Rename sheets
$wsCommands = $wb.Worksheets.Item("Sheet1")
$wsCommands.Name = "Commands"

$wsMisc = $wb.Worksheets.Item("Sheet2")
$wsMisc.Name = "Misc."

# Fix up the Print_Area defined name so it points at the renamed sheet
$printAreaName = $wb.Names.Item(1)
$printAreaName.RefersTo = "=Commands!`$A`$1:`$N`$58"

# Update the view on the "Commands" sheet: select cell C71
# (this also clears the old topLeftCell scroll-freeze from A25)
$wsCommands.Activate()
$wsCommands.Range("C71").Select()

# Update the view on the "Misc." sheet: select cell H21
$wsMisc.Activate()
$wsMisc.Range("H21").Select()

# Re-activate Commands sheet (tabSelected="1" in diff on Commands, consistent with original)
$wsCommands.Activate()
